# Updated scrapers & data for latest prices
#
# 1) "4x4 Squat Racks" sheet: cell E4 held a base64 data: URI placeholder
#    (leftover from a broken scrape) instead of a real image URL. Replace it
#    with the actual scraped image URL and wire it up as a clickable
#    hyperlink, matching the other "Image URL" / "Product Page" cells in the
#    table (blue, underlined).
# 2) "Squat Stands" sheet: row 5 ("No Squat Stand Available") used the
#    abbreviation "NA" for its Price and Product Page columns; spell it out
#    as "Not Available" for clarity.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "4x4 Squat Racks" ---
$ws1 = $wb.Worksheets.Item(1)

$newImageUrl = "https://garagegymlab.com/wp-content/uploads/Rogue-RM-3-Monster-Rack-2.0-Blue.jpg"

# Replace the stale base64 placeholder with the real image URL.
$ws1.Range("E4").Value = $newImageUrl

# Turn it into a real hyperlink (like every other Image URL / Product Page cell).
$ws1.Hyperlinks.Add($ws1.Range("E4"), $newImageUrl)

# Pick up the same visual style (blue + underline) the neighboring hyperlink
# cells already use, by copying formatting from F4 onto E4.
$ws1.Range("F4").Copy()
$ws1.Range("E4").PasteSpecial(-4122)

# --- Sheet 2: "Squat Stands" ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C5").Value = "Not Available"
$ws2.Range("F5").Value = "Not Available"
